$wb = $excel.ActiveWorkbook

# --- Rename existing headers to reflect PO terminology ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet as the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# --- Match page margins used by the other sheets (0.75/0.75/1/1/0.5/0.5 in) ---
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# --- Header row ---
$newSheet.Cells.Item(1,1).Value = "ds"
$newSheet.Cells.Item(1,2).Value = "PO_Forecast"
$newSheet.Cells.Item(1,3).Value = "yhat_lower"
$newSheet.Cells.Item(1,4).Value = "yhat_upper"

$headerRange = $newSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$newSheet.Cells.Item(2,1).Value = 44934.99999999999
$newSheet.Cells.Item(2,2).Value = 104
$newSheet.Cells.Item(2,3).Value = -264.5544824188984
$newSheet.Cells.Item(2,4).Value = 483.6636126322877
$newSheet.Cells.Item(3,1).Value = 44941.99999999999
$newSheet.Cells.Item(3,2).Value = 108
$newSheet.Cells.Item(3,3).Value = -279.2046031565401
$newSheet.Cells.Item(3,4).Value = 497.0945093418065
$newSheet.Cells.Item(4,1).Value = 44948.99999999999
$newSheet.Cells.Item(4,2).Value = 111
$newSheet.Cells.Item(4,3).Value = -267.7163719702024
$newSheet.Cells.Item(4,4).Value = 457.5212404555231
$newSheet.Cells.Item(5,1).Value = 44955.99999999999
$newSheet.Cells.Item(5,2).Value = 114
$newSheet.Cells.Item(5,3).Value = -273.6315711454195
$newSheet.Cells.Item(5,4).Value = 507.0358337451593
$newSheet.Cells.Item(6,1).Value = 44962.99999999999
$newSheet.Cells.Item(6,2).Value = 117
$newSheet.Cells.Item(6,3).Value = -250.081471906212
$newSheet.Cells.Item(6,4).Value = 498.2898967643828
$newSheet.Cells.Item(7,1).Value = 44969.99999999999
$newSheet.Cells.Item(7,2).Value = 120
$newSheet.Cells.Item(7,3).Value = -280.4855697488323
$newSheet.Cells.Item(7,4).Value = 497.4232307123951
$newSheet.Cells.Item(8,1).Value = 44976.99999999999
$newSheet.Cells.Item(8,2).Value = 124
$newSheet.Cells.Item(8,3).Value = -264.1442315628562
$newSheet.Cells.Item(8,4).Value = 499.1292316534843
$newSheet.Cells.Item(9,1).Value = 44983.99999999999
$newSheet.Cells.Item(9,2).Value = 127
$newSheet.Cells.Item(9,3).Value = -233.8494234056851
$newSheet.Cells.Item(9,4).Value = 498.4418865388262
$newSheet.Cells.Item(10,1).Value = 44990.99999999999
$newSheet.Cells.Item(10,2).Value = 130
$newSheet.Cells.Item(10,3).Value = -276.3185392744633
$newSheet.Cells.Item(10,4).Value = 508.8632200736048
$newSheet.Cells.Item(11,1).Value = 45032.99999999999
$newSheet.Cells.Item(11,2).Value = 149
$newSheet.Cells.Item(11,3).Value = -237.6555739995866
$newSheet.Cells.Item(11,4).Value = 530.9765280175613
$newSheet.Cells.Item(12,1).Value = 45081.99999999999
$newSheet.Cells.Item(12,2).Value = 172
$newSheet.Cells.Item(12,3).Value = -206.6104100792046
$newSheet.Cells.Item(12,4).Value = 535.4751292011185
$newSheet.Cells.Item(13,1).Value = 45088.99999999999
$newSheet.Cells.Item(13,2).Value = 175
$newSheet.Cells.Item(13,3).Value = -225.5481082656206
$newSheet.Cells.Item(13,4).Value = 551.383453190033
$newSheet.Cells.Item(14,1).Value = 45095.99999999999
$newSheet.Cells.Item(14,2).Value = 178
$newSheet.Cells.Item(14,3).Value = -189.0893159799625
$newSheet.Cells.Item(14,4).Value = 517.5927491491868
$newSheet.Cells.Item(15,1).Value = 45102.99999999999
$newSheet.Cells.Item(15,2).Value = 182
$newSheet.Cells.Item(15,3).Value = -197.6411589506774
$newSheet.Cells.Item(15,4).Value = 566.3607138802321
$newSheet.Cells.Item(16,1).Value = 45109.99999999999
$newSheet.Cells.Item(16,2).Value = 185
$newSheet.Cells.Item(16,3).Value = -185.1338025391354
$newSheet.Cells.Item(16,4).Value = 554.8849778438534
$newSheet.Cells.Item(17,1).Value = 45130.99999999999
$newSheet.Cells.Item(17,2).Value = 195
$newSheet.Cells.Item(17,3).Value = -193.0607688481891
$newSheet.Cells.Item(17,4).Value = 572.4984087915989
$newSheet.Cells.Item(18,1).Value = 45144.99999999999
$newSheet.Cells.Item(18,2).Value = 201
$newSheet.Cells.Item(18,3).Value = -153.020499731146
$newSheet.Cells.Item(18,4).Value = 551.6099433635899
$newSheet.Cells.Item(19,1).Value = 45165.99999999999
$newSheet.Cells.Item(19,2).Value = 211
$newSheet.Cells.Item(19,3).Value = -195.7957854077332
$newSheet.Cells.Item(19,4).Value = 597.8404683996341
$newSheet.Cells.Item(20,1).Value = 45172.99999999999
$newSheet.Cells.Item(20,2).Value = 214
$newSheet.Cells.Item(20,3).Value = -185.3408411905344
$newSheet.Cells.Item(20,4).Value = 610.8437420417745
$newSheet.Cells.Item(21,1).Value = 45179.99999999999
$newSheet.Cells.Item(21,2).Value = 217
$newSheet.Cells.Item(21,3).Value = -164.2686750123407
$newSheet.Cells.Item(21,4).Value = 625.2340319358416
$newSheet.Cells.Item(22,1).Value = 45186.99999999999
$newSheet.Cells.Item(22,2).Value = 220
$newSheet.Cells.Item(22,3).Value = -146.143448050621
$newSheet.Cells.Item(22,4).Value = 585.2039270884671
$newSheet.Cells.Item(23,1).Value = 45193.99999999999
$newSheet.Cells.Item(23,2).Value = 224
$newSheet.Cells.Item(23,3).Value = -168.3025319000546
$newSheet.Cells.Item(23,4).Value = 612.9226780078297
$newSheet.Cells.Item(24,1).Value = 45200.99999999999
$newSheet.Cells.Item(24,2).Value = 227
$newSheet.Cells.Item(24,3).Value = -137.53512931447
$newSheet.Cells.Item(24,4).Value = 612.4929487256844
$newSheet.Cells.Item(25,1).Value = 45207.99999999999
$newSheet.Cells.Item(25,2).Value = 230
$newSheet.Cells.Item(25,3).Value = -150.7093732048846
$newSheet.Cells.Item(25,4).Value = 642.3411946894041
$newSheet.Cells.Item(26,1).Value = 45214.99999999999
$newSheet.Cells.Item(26,2).Value = 233
$newSheet.Cells.Item(26,3).Value = -159.5230453904322
$newSheet.Cells.Item(26,4).Value = 606.7904009715575
$newSheet.Cells.Item(27,1).Value = 45221.99999999999
$newSheet.Cells.Item(27,2).Value = 236
$newSheet.Cells.Item(27,3).Value = -149.906347366687
$newSheet.Cells.Item(27,4).Value = 630.9229676828594
$newSheet.Cells.Item(28,1).Value = 45228.99999999999
$newSheet.Cells.Item(28,2).Value = 240
$newSheet.Cells.Item(28,3).Value = -163.4889411143204
$newSheet.Cells.Item(28,4).Value = 638.1307929982212
$newSheet.Cells.Item(29,1).Value = 45235.99999999999
$newSheet.Cells.Item(29,2).Value = 243
$newSheet.Cells.Item(29,3).Value = -148.9929232773535
$newSheet.Cells.Item(29,4).Value = 624.130745827126
$newSheet.Cells.Item(30,1).Value = 45242.99999999999
$newSheet.Cells.Item(30,2).Value = 246
$newSheet.Cells.Item(30,3).Value = -122.9672324196721
$newSheet.Cells.Item(30,4).Value = 646.3801668519885
$newSheet.Cells.Item(31,1).Value = 45249.99999999999
$newSheet.Cells.Item(31,2).Value = 249
$newSheet.Cells.Item(31,3).Value = -138.9421416870813
$newSheet.Cells.Item(31,4).Value = 627.6107153080986
$newSheet.Cells.Item(32,1).Value = 45270.99999999999
$newSheet.Cells.Item(32,2).Value = 259
$newSheet.Cells.Item(32,3).Value = -121.5845935359957
$newSheet.Cells.Item(32,4).Value = 604.2546800098048
$newSheet.Cells.Item(33,1).Value = 45277.99999999999
$newSheet.Cells.Item(33,2).Value = 262
$newSheet.Cells.Item(33,3).Value = -83.44479670365122
$newSheet.Cells.Item(33,4).Value = 626.0248953906053
$newSheet.Cells.Item(34,1).Value = 45298.99999999999
$newSheet.Cells.Item(34,2).Value = 272
$newSheet.Cells.Item(34,3).Value = -102.6154623575082
$newSheet.Cells.Item(34,4).Value = 676.1838616276142
$newSheet.Cells.Item(35,1).Value = 45305.99999999999
$newSheet.Cells.Item(35,2).Value = 275
$newSheet.Cells.Item(35,3).Value = -107.0288722903184
$newSheet.Cells.Item(35,4).Value = 657.3926345302947
$newSheet.Cells.Item(36,1).Value = 45319.99999999999
$newSheet.Cells.Item(36,2).Value = 282
$newSheet.Cells.Item(36,3).Value = -128.3904998434491
$newSheet.Cells.Item(36,4).Value = 648.6251569953067
$newSheet.Cells.Item(37,1).Value = 45326.99999999999
$newSheet.Cells.Item(37,2).Value = 285
$newSheet.Cells.Item(37,3).Value = -115.2794973991142
$newSheet.Cells.Item(37,4).Value = 642.8501503405352
$newSheet.Cells.Item(38,1).Value = 45333.99999999999
$newSheet.Cells.Item(38,2).Value = 288
$newSheet.Cells.Item(38,3).Value = -96.92855076088588
$newSheet.Cells.Item(38,4).Value = 680.6549013950346
$newSheet.Cells.Item(39,1).Value = 45340.99999999999
$newSheet.Cells.Item(39,2).Value = 291
$newSheet.Cells.Item(39,3).Value = -103.4564412217889
$newSheet.Cells.Item(39,4).Value = 664.5532189980509
$newSheet.Cells.Item(40,1).Value = 45347.99999999999
$newSheet.Cells.Item(40,2).Value = 294
$newSheet.Cells.Item(40,3).Value = -87.60498719282342
$newSheet.Cells.Item(40,4).Value = 677.8611181492525
$newSheet.Cells.Item(41,1).Value = 45354.99999999999
$newSheet.Cells.Item(41,2).Value = 298
$newSheet.Cells.Item(41,3).Value = -95.22803217145393
$newSheet.Cells.Item(41,4).Value = 706.7114531754858
$newSheet.Cells.Item(42,1).Value = 45361.99999999999
$newSheet.Cells.Item(42,2).Value = 301
$newSheet.Cells.Item(42,3).Value = -78.34314726230585
$newSheet.Cells.Item(42,4).Value = 697.1913714793385
$newSheet.Cells.Item(43,1).Value = 45368.99999999999
$newSheet.Cells.Item(43,2).Value = 304
$newSheet.Cells.Item(43,3).Value = -78.88770304403705
$newSheet.Cells.Item(43,4).Value = 665.4199392570455
$newSheet.Cells.Item(44,1).Value = 45375.99999999999
$newSheet.Cells.Item(44,2).Value = 307
$newSheet.Cells.Item(44,3).Value = -82.01046834497346
$newSheet.Cells.Item(44,4).Value = 691.6513198564716
$newSheet.Cells.Item(45,1).Value = 45382.99999999999
$newSheet.Cells.Item(45,2).Value = 311
$newSheet.Cells.Item(45,3).Value = -97.90383396352804
$newSheet.Cells.Item(45,4).Value = 689.3905529132196
$newSheet.Cells.Item(46,1).Value = 45389.99999999999
$newSheet.Cells.Item(46,2).Value = 314
$newSheet.Cells.Item(46,3).Value = -65.36472835647552
$newSheet.Cells.Item(46,4).Value = 702.6927715697251
$newSheet.Cells.Item(47,1).Value = 45403.99999999999
$newSheet.Cells.Item(47,2).Value = 320
$newSheet.Cells.Item(47,3).Value = -41.54195831201732
$newSheet.Cells.Item(47,4).Value = 711.6122291007453
$newSheet.Cells.Item(48,1).Value = 45410.99999999999
$newSheet.Cells.Item(48,2).Value = 323
$newSheet.Cells.Item(48,3).Value = -38.55179318204208
$newSheet.Cells.Item(48,4).Value = 731.2202042206683
$newSheet.Cells.Item(49,1).Value = 45417.99999999999
$newSheet.Cells.Item(49,2).Value = 327
$newSheet.Cells.Item(49,3).Value = -34.88782985927193
$newSheet.Cells.Item(49,4).Value = 699.2708759240123
$newSheet.Cells.Item(50,1).Value = 45424.99999999999
$newSheet.Cells.Item(50,2).Value = 330
$newSheet.Cells.Item(50,3).Value = -53.45992751878082
$newSheet.Cells.Item(50,4).Value = 712.2691656791239
$newSheet.Cells.Item(51,1).Value = 45431.99999999999
$newSheet.Cells.Item(51,2).Value = 333
$newSheet.Cells.Item(51,3).Value = -76.93123184600485
$newSheet.Cells.Item(51,4).Value = 738.1173262006205
$newSheet.Cells.Item(52,1).Value = 45438.99999999999
$newSheet.Cells.Item(52,2).Value = 336
$newSheet.Cells.Item(52,3).Value = -24.1006805992767
$newSheet.Cells.Item(52,4).Value = 734.9525482569555
$newSheet.Cells.Item(53,1).Value = 45445.99999999999
$newSheet.Cells.Item(53,2).Value = 340
$newSheet.Cells.Item(53,3).Value = -66.16986170477574
$newSheet.Cells.Item(53,4).Value = 724.9689794204226
$newSheet.Cells.Item(54,1).Value = 45452.99999999999
$newSheet.Cells.Item(54,2).Value = 343
$newSheet.Cells.Item(54,3).Value = -58.2199609311767
$newSheet.Cells.Item(54,4).Value = 737.0528665222686
$newSheet.Cells.Item(55,1).Value = 45459.99999999999
$newSheet.Cells.Item(55,2).Value = 346
$newSheet.Cells.Item(55,3).Value = -28.32786583160408
$newSheet.Cells.Item(55,4).Value = 733.5483479019557
$newSheet.Cells.Item(56,1).Value = 45473.99999999999
$newSheet.Cells.Item(56,2).Value = 352
$newSheet.Cells.Item(56,3).Value = -41.86272858779922
$newSheet.Cells.Item(56,4).Value = 742.4894471824673
$newSheet.Cells.Item(57,1).Value = 45487.99999999999
$newSheet.Cells.Item(57,2).Value = 359
$newSheet.Cells.Item(57,3).Value = -18.7859909752304
$newSheet.Cells.Item(57,4).Value = 741.2234271517837
$newSheet.Cells.Item(58,1).Value = 45578.99999999999
$newSheet.Cells.Item(58,2).Value = 401
$newSheet.Cells.Item(58,3).Value = 6.679172267878928
$newSheet.Cells.Item(58,4).Value = 760.8853506168698
$newSheet.Cells.Item(59,1).Value = 45585.99999999999
$newSheet.Cells.Item(59,2).Value = 404
$newSheet.Cells.Item(59,3).Value = 3.433391083362421
$newSheet.Cells.Item(59,4).Value = 772.6210741277872
$newSheet.Cells.Item(60,1).Value = 45592.99999999999
$newSheet.Cells.Item(60,2).Value = 407
$newSheet.Cells.Item(60,3).Value = 38.95075191517272
$newSheet.Cells.Item(60,4).Value = 770.2470551871862
$newSheet.Cells.Item(61,1).Value = 45599.99999999999
$newSheet.Cells.Item(61,2).Value = 410
$newSheet.Cells.Item(61,3).Value = 23.3076015071378
$newSheet.Cells.Item(61,4).Value = 801.7385100048101
$newSheet.Cells.Item(62,1).Value = 45613.99999999999
$newSheet.Cells.Item(62,2).Value = 417
$newSheet.Cells.Item(62,3).Value = 45.19436953447708
$newSheet.Cells.Item(62,4).Value = 815.381681777607
$newSheet.Cells.Item(63,1).Value = 45620.99999999999
$newSheet.Cells.Item(63,2).Value = 420
$newSheet.Cells.Item(63,3).Value = 51.20812371728365
$newSheet.Cells.Item(63,4).Value = 793.4445121088487
$newSheet.Cells.Item(64,1).Value = 45627.99999999999
$newSheet.Cells.Item(64,2).Value = 423
$newSheet.Cells.Item(64,3).Value = -1.027017586358086
$newSheet.Cells.Item(64,4).Value = 813.9460319728528
$newSheet.Cells.Item(65,1).Value = 45634.99999999999
$newSheet.Cells.Item(65,2).Value = 427
$newSheet.Cells.Item(65,3).Value = 44.16784170705583
$newSheet.Cells.Item(65,4).Value = 768.9661968461363
$newSheet.Cells.Item(66,1).Value = 45641.99999999999
$newSheet.Cells.Item(66,2).Value = 430
$newSheet.Cells.Item(66,3).Value = 44.16850666335547
$newSheet.Cells.Item(66,4).Value = 803.6171143420129
$newSheet.Cells.Item(67,1).Value = 45648.99999999999
$newSheet.Cells.Item(67,2).Value = 433
$newSheet.Cells.Item(67,3).Value = 43.19269030582247
$newSheet.Cells.Item(67,4).Value = 833.6594584809619
$newSheet.Cells.Item(68,1).Value = 45655.99999999999
$newSheet.Cells.Item(68,2).Value = 436
$newSheet.Cells.Item(68,3).Value = 70.65510752189404
$newSheet.Cells.Item(68,4).Value = 826.0595223098703
$newSheet.Cells.Item(69,1).Value = 45662.99999999999
$newSheet.Cells.Item(69,2).Value = 439
$newSheet.Cells.Item(69,3).Value = 49.72245673124184
$newSheet.Cells.Item(69,4).Value = 852.8903620465629
$newSheet.Cells.Item(70,1).Value = 45669.99999999999
$newSheet.Cells.Item(70,2).Value = 443
$newSheet.Cells.Item(70,3).Value = 50.32779553698767
$newSheet.Cells.Item(70,4).Value = 830.9198299343171

# --- Date/time number format for column A data rows (matches other sheets) ---
$newSheet.Range("A2:A70").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "PO Forecast sheet created with $($newSheet.UsedRange.Rows.Count) rows"
